$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("batteries")

# Rename the header in Q1 from "default_operation_mode" to "operation_mode"
$ws.Range("Q1").Value = "operation_mode"

# Update the active cell selection to K11, matching the saved view state
$ws.Range("K11").Select()
